$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = "2023-12-06 16:18:42"
$ws.Range("B28").Value = 0.001

$ws.Range("A29").Value = "2023-12-06 16:19:57"
$ws.Range("B29").Value = 0.005000000000000001

$ws.Range("A30").Value = "2023-12-06 16:20:45"
$ws.Range("B30").Value = 0.0022
